# Applies the Sep 25 2024 cryptos price/volume refresh to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.329.99'
$ws.Range('E2').Value = '  -0.35%  '
# Row 3
$ws.Range('D3').Value = '2.593.87'
$ws.Range('E3').Value = '  -1.66%  '
# Row 4
$ws.Range('E4').Value = '  +0.17%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.16'
$ws.Range('E5').Value = '  -2.88%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.40'
$ws.Range('E6').Value = '  +1.90%  '
# Row 7
$ws.Range('E7').Value = '  +0.13%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.585'
$ws.Range('E8').Value = '  -0.20%  '
# Row 9
$ws.Range('E9').Value = '  +0.97%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.65'
$ws.Range('E10').Value = '  +2.04%  '
# Row 11
$ws.Range('E11').Value = '  +0.84%  '
# Row 12
$ws.Range('E12').Value = '  -0.76%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.41'
$ws.Range('E13').Value = '  +0.00%  '
# Row 14
$ws.Range('D14').Value = '3.058.50'
$ws.Range('E14').Value = '  -1.45%  '
# Row 15
$ws.Range('D15').Value = '63.216.09'
$ws.Range('E15').Value = '  -0.25%  '
# Row 16
$ws.Range('E16').Value = '  +6.11%  '
# Row 17
$ws.Range('D17').Value = '2.596.22'
$ws.Range('E17').Value = '  -1.16%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.31'
$ws.Range('E18').Value = '  +5.45%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.71'
$ws.Range('E19').Value = '  +3.34%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.34'
$ws.Range('E20').Value = '  +0.11%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.88'
$ws.Range('E21').Value = '  -0.32%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.17%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.58'
$ws.Range('E23').Value = '  +2.11%  '
# Row 24
$ws.Range('E24').Value = '  +3.44%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.21'
$ws.Range('E25').Value = '  +1.87%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.66'
$ws.Range('E26').Value = '  -1.69%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '552.48'
$ws.Range('E27').Value = '  -0.65%  '
# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.01'
$ws.Range('E28').Value = '  +0.70%  '
# Row 29
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.97'
$ws.Range('E29').Value = '  -1.23%  '
# Row 30
$ws.Range('E30').Value = '  -1.14%  '
# Row 31
$ws.Range('E31').Value = '  +0.67%  '
# Row 32
$ws.Range('D32').Value = '0.0₃0843'
$ws.Range('E32').Value = '  -0.76%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.74'
$ws.Range('E33').Value = '  -1.21%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.17'
$ws.Range('E34').Value = '  -3.30%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '167.25'
$ws.Range('E35').Value = '  -0.96%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.412'
$ws.Range('E36').Value = '  +1.96%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.25%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.54'
$ws.Range('E38').Value = '  +2.48%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.92'
$ws.Range('E39').Value = '  +0.02%  '
# Row 40
$ws.Range('E40').Value = '  -0.02%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '166.77'
$ws.Range('E41').Value = '  +1.11%  '
# Row 42
$ws.Range('E42').Value = '  -0.71%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.94'
$ws.Range('E43').Value = '  +4.27%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0581'
$ws.Range('E44').Value = '  +2.35%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.33'
$ws.Range('E45').Value = '  +1.44%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.630'
$ws.Range('E46').Value = '  +0.35%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0251'
$ws.Range('E47').Value = '  +2.71%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.02'
$ws.Range('E48').Value = '  +1.87%  '
# Row 49
$ws.Range('E49').Value = '  +0.81%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.07'
$ws.Range('E50').Value = '  +1.50%  '
# Row 51
$ws.Range('E51').Value = '  +18.37%  '
